# Edits to the guide
# Adjusts the positions/sizes of several shapes on slide 1 of the
# FortiGate architecture diagram and expands the "(main)"/"(host)"
# labels into "(main, BYOL)" / "(host, on-demand)".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Id -eq $id) {
            return $sh
        }
    }
    return $null
}

# EMU -> point helper (PowerPoint COM works in points; OOXML stores EMU).
#
# The host's Shape.Left/Top/Width/Height setters convert the incoming
# point value to EMU as floor(float32(points) * 12700) — i.e. the double
# we pass gets narrowed to a 32-bit float before the EMU conversion, and
# the EMU is then floored rather than rounded. A naive `emu / 12700.0`
# point value can therefore land the float32 product just below the
# intended integer EMU and save one EMU short. To hit an exact target
# EMU we search a small neighbourhood of doubles around the naive
# quotient until the same float32 -> floor(*12700) pipeline reproduces
# the exact target.
function Emu([double]$targetEmu) {
    $base = $targetEmu / 12700.0
    $step = 0.0000001
    for ($k = 0; $k -le 4000; $k++) {
        foreach ($sign in @(1, -1)) {
            if ($k -eq 0 -and $sign -eq -1) { continue }
            $cand = $base + ($sign * $k * $step)
            $f = [float]$cand
            $emu = [Math]::Floor([double]$f * 12700.0)
            if ($emu -eq $targetEmu) {
                return $cand
            }
        }
    }
    # Fallback (should not happen for the values used in this script).
    return $base
}

# --- Shape id=59 "FortiGate" label (near top-right FortiGate icon) ---
$sh59 = Get-ShapeById $s 59
$sh59.Left = Emu 2895966
$sh59.Top  = Emu 1944284

# --- Shape id=71 "Auto Scaling group" dashed box (top) ---
$sh71 = Get-ShapeById $s 71
$sh71.Left   = Emu 1264436
$sh71.Top    = Emu 1603778
$sh71.Width  = Emu 2221535
$sh71.Height = Emu 602068

# --- Shape id=85 "Auto Scaling group" dashed box (bottom) ---
$sh85 = Get-ShapeById $s 85
$sh85.Left   = Emu 808518
$sh85.Top    = Emu 2381802
$sh85.Width  = Emu 3120695
$sh85.Height = Emu 608230

# --- Shape id=88 small picture near top FortiGate ---
$sh88 = Get-ShapeById $s 88
$sh88.Left = Emu 2308446
$sh88.Top  = Emu 1597004

# --- Shape id=100 "FortiGate (main)" label ---
$sh100 = Get-ShapeById $s 100
$sh100.Left = Emu 1213076
$sh100.Top  = Emu 1951272
$tr100 = $sh100.TextFrame.TextRange
$para2 = $tr100.Paragraphs(2)
# Paragraph text is "(main)"; replace "main)" (chars 2-6) with "main, BYOL)"
# so the split lands after the opening "(", producing two runs: "(" and
# "main, BYOL)" with identical run formatting (matches target OOXML).
$rng100 = $para2.Characters(2, 5)
$rng100.Text = "main, BYOL)"

# --- Shape id=102 "FortiGate<br/>(host)" label ---
$sh102 = Get-ShapeById $s 102
$sh102.Left   = Emu 808518
$sh102.Top    = Emu 2728441
$sh102.Width  = Emu 791682
$tr102 = $sh102.TextFrame.TextRange
# Full text is "FortiGate" + line break + "(host)" (line break counts as
# one character), so "(host)" starts at character 11; replace "host)"
# (chars 12-16) with "host, on-demand)" to split into "(" + "host, on-demand)".
$rng102 = $tr102.Characters(12, 5)
$rng102.Text = "host, on-demand)"

# --- Shape id=110 connector (right, near top FortiGate) ---
$sh110 = Get-ShapeById $s 110
$sh110.Left = Emu 3220439

# --- Shape id=111 connector (left, near top FortiGate) ---
$sh111 = Get-ShapeById $s 111
$sh111.Left = Emu 1539867

# --- Shape id=119 small "Graphic 56" picture (right) ---
$sh119 = Get-ShapeById $s 119
$sh119.Left = Emu 3101854

# --- Shape id=120 small "Graphic 56" picture (left) ---
$sh120 = Get-ShapeById $s 120
$sh120.Left = Emu 1421829

# --- Shape id=52 small picture near bottom FortiGate ---
$sh52 = Get-ShapeById $s 52
$sh52.Left   = Emu 2308446
$sh52.Top    = Emu 2374523
$sh52.Width  = Emu 136534
$sh52.Height = Emu 141685
